# Update cryptocurrency price and volume(1h) data as of Mon Dec  4 10:52:39 UTC 2023
# (GitHub Actions scheduled refresh of cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.035.34'
$ws.Range('E2').Value = '  +6.32%  '
$ws.Range('D3').Value = '2.271.17'
$ws.Range('E3').Value = '  +4.76%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '234.82'
$ws.Range('E5').Value = '  +2.58%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.645'
$ws.Range('E6').Value = '  +1.62%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '63.63'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E9').Value = '  +3.98%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '60.15'
$ws.Range('E10').Value = '  +3.55%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0898'
$ws.Range('E11').Value = '  +5.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.106'
$ws.Range('E12').Value = '  +1.90%  '
$ws.Range('D13').Value = '2.602.65'
$ws.Range('E13').Value = '  +4.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '16.12'
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '22.98'
$ws.Range('E15').Value = '  +4.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.824'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('E17').Value = '  +3.54%  '
$ws.Range('D18').Value = '2.265.60'
$ws.Range('E18').Value = '  +4.59%  '
$ws.Range('D19').Value = '41.855.13'
$ws.Range('E19').Value = '  +5.97%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '74.83'
$ws.Range('E20').Value = '  +3.85%  '
$ws.Range('D21').Value = '0.0₃0933'
$ws.Range('E21').Value = '  +9.91%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.17'
$ws.Range('E22').Value = '  -0.68%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '253.03'
$ws.Range('E23').Value = '  +10.14%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('E25').Value = '  +3.46%  '
$ws.Range('E26').Value = '  +2.71%  '
$ws.Range('E27').Value = '  +7.57%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.87'
$ws.Range('E28').Value = '  +3.13%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '171.27'
$ws.Range('E29').Value = '  -0.65%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '20.54'
$ws.Range('E30').Value = '  +3.20%  '
$ws.Range('E31').Value = '  +2.27%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.83'
$ws.Range('E32').Value = '  +7.04%  '
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.11'
$ws.Range('E34').Value = '  +7.88%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.81'
$ws.Range('E35').Value = '  +4.15%  '
$ws.Range('E36').Value = '  +3.31%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.85'
$ws.Range('E37').Value = '  -3.02%  '
$ws.Range('E38').Value = '  +6.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.46'
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.000262'
$ws.Range('E40').Value = '  +50.58%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.16'
$ws.Range('E41').Value = '  +18.68%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('E43').Value = '  +5.88%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.69'
$ws.Range('E44').Value = '  +9.68%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '102.22'
$ws.Range('E45').Value = '  -0.64%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '17.77'
$ws.Range('E46').Value = '  -1.15%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.24'
$ws.Range('E47').Value = '  +3.27%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0987'
$ws.Range('E48').Value = '  +6.62%  '
$ws.Range('D49').Value = '1.506.52'
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('E50').Value = '  +1.69%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.82'
$ws.Range('E51').Value = '  +0.03%  '
